# refs #554 Code Review dokumentiert
#
# The commit merges runs that had been split apart purely for the
# benefit of the spell-checker (w:proofErr spellStart/spellEnd wrapping
# individual "unknown" words like Gfeller, Treichler, VisualStateGroups,
# etc). Re-typing/replacing the full span of text that used to be spread
# across multiple runs collapses it back into a single run and drops the
# now orphaned w:proofErr markers, matching the target XML.
#
# It also relocates the "_GoBack" bookmark from just before the
# "16.12.2011" heading to the end of the preceding list item (right
# after "...Dictionaries)."), removing the trailing run that used to
# hold a single space character there.

$d = $word.ActiveDocument
$wdFindContinue = 1
$wdReplaceNone = 0
$wdCollapseEnd = 0

function Merge-Text($oldText, $newText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, 2)
    if (-not $ok) {
        throw "Find failed for: $oldText"
    }
}

# 1. "Zusätzlich ... Michael " + "Gfeller" + " durchgeführt ..."
Merge-Text `
    "Zusätzlich zu den internen Code Reviews wurden auch noch Code Reviews mit Christian Moser und Michael Gfeller durchgeführt. Die Details dazu sind untenstehend aufgelistet." `
    "Zusätzlich zu den internen Code Reviews wurden auch noch Code Reviews mit Christian Moser und Michael Gfeller durchgeführt. Die Details dazu sind untenstehend aufgelistet."

# 2. "Anwesend: Christian Moser, Christina Heidt, Delia " + "Treichler" + ", Lukas Elmer"
Merge-Text `
    "Anwesend: Christian Moser, Christina Heidt, Delia Treichler, Lukas Elmer" `
    "Anwesend: Christian Moser, Christina Heidt, Delia Treichler, Lukas Elmer"

# 3. "Um die Animationen ..." + "VisualStateGroups" + " ... Expression " + "Blend" + " "
Merge-Text `
    "Um die Animationen zu gestalten, wurden VisualStateGroups eingesetzt. Diese wurden im Team besprochen und im Expression Blend " `
    "Um die Animationen zu gestalten, wurden VisualStateGroups eingesetzt. Diese wurden im Team besprochen und im Expression Blend "

# 4. "Review mit: Michael " + "Gfeller" (first occurrence, 9.12.2011 section)
Merge-Text `
    "Review mit: Michael Gfeller" `
    "Review mit: Michael Gfeller"

# 5. "Anwesend: Michael " + "Gfeller" + ", Christina Heidt, Delia " + "Treichler" + ", Lukas Elmer"
Merge-Text `
    "Anwesend: Michael Gfeller, Christina Heidt, Delia Treichler, Lukas Elmer" `
    "Anwesend: Michael Gfeller, Christina Heidt, Delia Treichler, Lukas Elmer"

# 6. Big "Da das XAML des ..." run pile-up, merged down to two runs.
Merge-Text `
    "Da das XAML des OverviewWindows zu lang und gross wurde, sollen die Styles in ein Styles.xaml ausgelagert werden. Diese können dann mithilfe eines ResourceDirectory inds Xaml eingebunden werden " `
    "Da das XAML des OverviewWindows zu lang und gross wurde, sollen die Styles in ein Styles.xaml ausgelagert werden. Diese können dann mithilfe eines ResourceDirectory inds Xaml eingebunden werden "

Merge-Text `
    "(Merged Resource Dictionaries)." `
    "(Merged Resource Dictionaries)."

# Drop the trailing single-space run after "...Dictionaries)." and move
# the _GoBack bookmark to sit right at the end of that paragraph.
Merge-Text "Dictionaries). " "Dictionaries)."

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()

$rng = $d.Content
$rng.Find.Execute("Dictionaries).") | Out-Null
$rng.Collapse($wdCollapseEnd)
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null

# 7. "nicht " + "freezed" + " "
Merge-Text `
    "nicht freezed " `
    "nicht freezed "

# 8. " Memory " + "Leaks"
Merge-Text `
    " Memory Leaks" `
    " Memory Leaks"

# 9. "ScrollToTopBehavior" + ": " + "DependencyPropertyDescriptor" + " ist statisches Konstrukt (" + "ItemsSourceProperty" + "): prüfen, dass es sich abmeldet"
Merge-Text `
    "ScrollToTopBehavior: DependencyPropertyDescriptor ist statisches Konstrukt (ItemsSourceProperty): prüfen, dass es sich abmeldet" `
    "ScrollToTopBehavior: DependencyPropertyDescriptor ist statisches Konstrukt (ItemsSourceProperty): prüfen, dass es sich abmeldet"

# ", sonst könnte hier ein Memory " + "Leak" + " entstehen."
Merge-Text `
    ", sonst könnte hier ein Memory Leak entstehen." `
    ", sonst könnte hier ein Memory Leak entstehen."

# 10. ", " + "Konstruktor" + ", dann private"
Merge-Text `
    ", Konstruktor, dann private" `
    ", Konstruktor, dann private"

# 11. " soll in ein " + "Konfigurations" + " File ausgelagert werden."
Merge-Text `
    " soll in ein Konfigurations File ausgelagert werden." `
    " soll in ein Konfigurations File ausgelagert werden."

# 12. "Surface" + " 2 zu benutzen."
Merge-Text `
    "Surface 2 zu benutzen." `
    "Surface 2 zu benutzen."

# 13. "Code dokumentieren für " + "public"
Merge-Text `
    "Code dokumentieren für public" `
    "Code dokumentieren für public"

# 14. "wenn nicht klar ist, um was es sich handelt, z.B. bei " + "Preload" + "()"
Merge-Text `
    "wenn nicht klar ist, um was es sich handelt, z.B. bei Preload()" `
    "wenn nicht klar ist, um was es sich handelt, z.B. bei Preload()"

# " im " + "ProjectNote" + " Model"
Merge-Text `
    " im ProjectNote Model" `
    " im ProjectNote Model"

# 16. "Review " + "mit" + ": Michael " + "Gfeller" (second occurrence, 16.12.2011 section)
Merge-Text `
    "Review mit: Michael Gfeller" `
    "Review mit: Michael Gfeller"

Write-Output "done"
